$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2508.625
$ws.Range("I4").Value = 2508.625
$ws.Range("K4").Value = 2508.625
$ws.Range("M4").Value = -2394.625
$ws.Range("H18").Value = 7595.0713
$ws.Range("I18").Value = 486.84616
$ws.Range("K18").Value = 486.84616
$ws.Range("M18").Value = -202.84616
$ws.Range("H129").Value = 309516.44
$ws.Range("I129").Value = 638.6
$ws.Range("J129").Value = 348126.2
$ws.Range("K129").Value = 1915.8
$ws.Range("L129").Value = 1044378.6
$ws.Range("M129").Value = 3084.2
$ws.Range("N129").Value = -1054378.6
$ws.Range("H137").Value = 2112.8462
$ws.Range("I137").Value = 1577.4445
$ws.Range("J137").Value = 3317.5
$ws.Range("K137").Value = 4732.333500000001
$ws.Range("L137").Value = 9952.5
$ws.Range("M137").Value = -2182.333500000001
$ws.Range("N137").Value = -15052.5
$ws.Range("H138").Value = 3602.861
$ws.Range("I138").Value = 1250.8667
$ws.Range("J138").Value = 5282.857
$ws.Range("K138").Value = 3752.6001
$ws.Range("L138").Value = 15848.571
$ws.Range("M138").Value = 1387.3999
$ws.Range("N138").Value = -26128.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 46875.09
$ws.Range("I2").Value = 1304.2142
$ws.Range("K2").Value = 1304.2142
$ws.Range("M2").Value = -1191.2142
$ws.Range("H45").Value = 1760.6923
$ws.Range("I45").Value = 1491.25
$ws.Range("J45").Value = 2191.8
$ws.Range("K45").Value = 1491.25
$ws.Range("L45").Value = 2191.8
$ws.Range("M45").Value = -1114.25
$ws.Range("N45").Value = -2945.8
$ws.Range("H107").Value = 27400.5
$ws.Range("J107").Value = 27400.5
$ws.Range("L107").Value = 27400.5
$ws.Range("N107").Value = -35080.5
$ws.Range("H116").Value = 46875.09
$ws.Range("I116").Value = 1304.2142
$ws.Range("K116").Value = 1304.2142
$ws.Range("M116").Value = 989.7858000000001
$ws.Range("H124").Value = 23929.857
$ws.Range("J124").Value = 23929.857
$ws.Range("L124").Value = 23929.857
$ws.Range("N124").Value = -33749.857
$ws.Range("H132").Value = 3781.3333
$ws.Range("I132").Value = 3781.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11343.9999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8813.999899999999
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 46875.09
$ws.Range("I3").Value = 1304.2142
$ws.Range("K3").Value = 1304.2142
$ws.Range("M3").Value = -1190.2142
$ws.Range("H105").Value = 119137.65
$ws.Range("I105").Value = 92357.27
$ws.Range("J105").Value = 168235
$ws.Range("K105").Value = 92357.27
$ws.Range("L105").Value = 168235
$ws.Range("M105").Value = -90610.27
$ws.Range("N105").Value = -171729

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38666.5
$ws.Range("I31").Value = 68195.664
$ws.Range("J31").Value = 4594.385
$ws.Range("K31").Value = 68195.664
$ws.Range("L31").Value = 4594.385
$ws.Range("M31").Value = -67900.664
$ws.Range("N31").Value = -5184.385
$ws.Range("H34").Value = 38666.5
$ws.Range("I34").Value = 68195.664
$ws.Range("J34").Value = 4594.385
$ws.Range("K34").Value = 68195.664
$ws.Range("L34").Value = 4594.385
$ws.Range("M34").Value = -67993.664
$ws.Range("N34").Value = -4998.385
$ws.Range("H50").Value = 14285
$ws.Range("J50").Value = 14285
$ws.Range("L50").Value = 14285
$ws.Range("N50").Value = -15535
$ws.Range("H51").Value = 7932.4165
$ws.Range("J51").Value = 7918.091
$ws.Range("L51").Value = 7918.091
$ws.Range("N51").Value = -9390.091
$ws.Range("H60").Value = 11067.5
$ws.Range("J60").Value = 11067.5
$ws.Range("L60").Value = 11067.5
$ws.Range("N60").Value = -12089.5
$ws.Range("H61").Value = 7932.4165
$ws.Range("J61").Value = 7918.091
$ws.Range("L61").Value = 7918.091
$ws.Range("N61").Value = -8614.091
$ws.Range("H68").Value = 14356.429
$ws.Range("J68").Value = 14356.429
$ws.Range("L68").Value = 14356.429
$ws.Range("N68").Value = -15854.429
$ws.Range("H71").Value = 14356.429
$ws.Range("J71").Value = 14356.429
$ws.Range("L71").Value = 43069.287
$ws.Range("N71").Value = -50557.287
$ws.Range("H86").Value = 2017.0714
$ws.Range("I86").Value = 1875
$ws.Range("K86").Value = 1875
$ws.Range("M86").Value = -752
$ws.Range("H89").Value = 2017.0714
$ws.Range("I89").Value = 1875
$ws.Range("K89").Value = 9375
$ws.Range("M89").Value = -3759
$ws.Range("H124").Value = 36994
$ws.Range("J124").Value = 36994
$ws.Range("L124").Value = 36994
$ws.Range("N124").Value = -41904
$ws.Range("H134").Value = 1380.4482
$ws.Range("I134").Value = 1309.5385
$ws.Range("J134").Value = 1995
$ws.Range("K134").Value = 3928.6155
$ws.Range("L134").Value = 5985
$ws.Range("M134").Value = -1393.6155
$ws.Range("N134").Value = -11055

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6373.2896
$ws.Range("I5").Value = 1070
$ws.Range("J5").Value = 15464.643
$ws.Range("K5").Value = 3210
$ws.Range("L5").Value = 46393.929
$ws.Range("M5").Value = -3098
$ws.Range("N5").Value = -46617.929
$ws.Range("H37").Value = 689872.2
$ws.Range("J37").Value = 689872.2
$ws.Range("L37").Value = 2069616.6
$ws.Range("N37").Value = -2069840.6
$ws.Range("H135").Value = 6373.2896
$ws.Range("I135").Value = 1070
$ws.Range("J135").Value = 15464.643
$ws.Range("K135").Value = 9630
$ws.Range("L135").Value = 139181.787
$ws.Range("M135").Value = -7095
$ws.Range("N135").Value = -144251.787

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 738.53845
$ws.Range("I22").Value = 499
$ws.Range("J22").Value = 758.5
$ws.Range("K22").Value = 499
$ws.Range("L22").Value = 758.5
$ws.Range("M22").Value = -204
$ws.Range("N22").Value = -1348.5
$ws.Range("H27").Value = 738.53845
$ws.Range("I27").Value = 499
$ws.Range("J27").Value = 758.5
$ws.Range("K27").Value = 499
$ws.Range("L27").Value = 758.5
$ws.Range("M27").Value = -392
$ws.Range("N27").Value = -972.5
$ws.Range("H122").Value = 1641.75
$ws.Range("I122").Value = 1386.1428
$ws.Range("J122").Value = 1999.6
$ws.Range("K122").Value = 4158.428400000001
$ws.Range("L122").Value = 5998.799999999999
$ws.Range("M122").Value = -1708.428400000001
$ws.Range("N122").Value = -10898.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 18400
$ws.Range("J64").Value = 18400
$ws.Range("L64").Value = 18400
$ws.Range("N64").Value = -18896
$ws.Range("H67").Value = 18400
$ws.Range("J67").Value = 18400
$ws.Range("L67").Value = 18400
$ws.Range("N67").Value = -20116
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H81").Value = 401135.8
$ws.Range("I81").Value = 1000000
$ws.Range("J81").Value = 251419.75
$ws.Range("K81").Value = 2000000
$ws.Range("L81").Value = 502839.5
$ws.Range("M81").Value = -1998939
$ws.Range("N81").Value = -504961.5
$ws.Range("H82").Value = 41980
$ws.Range("J82").Value = 41980
$ws.Range("L82").Value = 41980
$ws.Range("N82").Value = -42746
$ws.Range("H84").Value = 401135.8
$ws.Range("I84").Value = 1000000
$ws.Range("J84").Value = 251419.75
$ws.Range("K84").Value = 10000000
$ws.Range("L84").Value = 2514197.5
$ws.Range("M84").Value = -9994696
$ws.Range("N84").Value = -2524805.5
$ws.Range("H85").Value = 41980
$ws.Range("J85").Value = 41980
$ws.Range("L85").Value = 41980
$ws.Range("N85").Value = -44632
$ws.Range("H94").Value = 16025
$ws.Range("J94").Value = 16025
$ws.Range("L94").Value = 16025
$ws.Range("N94").Value = -17827
$ws.Range("H107").Value = 258000
$ws.Range("I107").Value = 93333.336
$ws.Range("K107").Value = 280000.008
$ws.Range("M107").Value = -278080.008
$ws.Range("H126").Value = 1470.3478
$ws.Range("I126").Value = 1467.5238
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 4402.5714
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -1932.5714
$ws.Range("N126").Value = -9440
$ws.Range("H136").Value = 1357.5143
$ws.Range("I136").Value = 522.3043
$ws.Range("J136").Value = 2958.3333
$ws.Range("K136").Value = 1566.9129
$ws.Range("L136").Value = 8874.999899999999
$ws.Range("M136").Value = 983.0871
$ws.Range("N136").Value = -13974.9999
